$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.391.70'
$ws.Range("E2").Value = '  +0.48%  '

$ws.Range("D3").Value = '2.553.40'
$ws.Range("E3").Value = '  -2.26%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.02'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.88%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.16%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.531'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.36%  '

$ws.Range("D9").Value = '2.552.49'
$ws.Range("E9").Value = '  -2.24%  '

$ws.Range("E10").Value = '  +1.24%  '

$ws.Range("E12").Value = '  -2.85%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.12'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.37%  '

$ws.Range("D15").Value = '3.010.86'
$ws.Range("E15").Value = '  -2.54%  '

$ws.Range("E16").Value = '  +0.04%  '

$ws.Range("D17").Value = '67.281.25'
$ws.Range("E17").Value = '  +0.37%  '

$ws.Range("D18").Value = '2.552.58'
$ws.Range("E18").Value = '  -1.00%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.19%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.43'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '356.18'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.33%  '

$ws.Range("E22").Value = '  -1.01%  '

$ws.Range("E23").Value = '  +1.42%  '

$ws.Range("E24").Value = '  +6.81%  '

$ws.Range("E25").Value = '  +0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.28'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.61%  '

$ws.Range("E27").Value = '  -3.39%  '

$ws.Range("D28").Value = '2.686.93'

$ws.Range("E29").Value = '  +0.10%  '

$ws.Range("E30").Value = '  +0.89%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '535.19'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.33%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.27'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.37%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.39%  '

$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("E35").Value = '  +0.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.00'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.09%  '

$ws.Range("E37").Value = '  +0.82%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.36'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("E39").Value = '  -0.29%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.32%  '

$ws.Range("E41").Value = '  -1.52%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.62%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.22'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.03%  '

$ws.Range("E44").Value = '  +6.47%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.79'
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '151.38'
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.566'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("E49").Value = '  -4.58%  '

$ws.Range("E50").Value = '  -0.84%  '

$ws.Range("E51").Value = '  +1.75%  '
